$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 1017, pushing the existing data
# (rows 1017-1104) down to rows 1019-1106. Excel's row-insert copies
# the formatting (e.g. the date style on column D) from the row above,
# matching the existing sheet's formatting for this block of data.
$ws.Rows.Item(1017).Insert()
$ws.Rows.Item(1017).Insert()

# New weekly data point for "Primera" quality, inserted at row 1017.
$ws.Cells.Item(1017, 1).Value2 = 6
$ws.Cells.Item(1017, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1017, 3).Value2 = "Metropolitana"
$ws.Cells.Item(1017, 4).Value2 = 45106
$ws.Cells.Item(1017, 5).Value2 = 13
$ws.Cells.Item(1017, 6).Value2 = 100112017
$ws.Cells.Item(1017, 7).Value2 = "Apio"
$ws.Cells.Item(1017, 8).Value2 = "Americana (o)"
$ws.Cells.Item(1017, 9).Value2 = "Primera"
$ws.Cells.Item(1017, 10).Value2 = 1200
$ws.Cells.Item(1017, 11).Value2 = 6000
$ws.Cells.Item(1017, 12).Value2 = 7000
$ws.Cells.Item(1017, 13).Value2 = 6458
$ws.Cells.Item(1017, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(1017, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(1017, 16).Value2 = 1076
$ws.Cells.Item(1017, 17).Value2 = 6
$ws.Cells.Item(1017, 18).Value2 = "Hortaliza"

# New weekly data point for "Segunda" quality, inserted at row 1018.
$ws.Cells.Item(1018, 1).Value2 = 6
$ws.Cells.Item(1018, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1018, 3).Value2 = "Metropolitana"
$ws.Cells.Item(1018, 4).Value2 = 45106
$ws.Cells.Item(1018, 5).Value2 = 13
$ws.Cells.Item(1018, 6).Value2 = 100112017
$ws.Cells.Item(1018, 7).Value2 = "Apio"
$ws.Cells.Item(1018, 8).Value2 = "Americana (o)"
$ws.Cells.Item(1018, 9).Value2 = "Segunda"
$ws.Cells.Item(1018, 10).Value2 = 650
$ws.Cells.Item(1018, 11).Value2 = 4000
$ws.Cells.Item(1018, 12).Value2 = 5000
$ws.Cells.Item(1018, 13).Value2 = 4538
$ws.Cells.Item(1018, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(1018, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(1018, 16).Value2 = 756
$ws.Cells.Item(1018, 17).Value2 = 6
$ws.Cells.Item(1018, 18).Value2 = "Hortaliza"
